# Weekly update: insert two new "current week" price rows for Betarraga at
# Femacal de La Calera, pushing the existing historical rows down by two
# positions (rows 266-384 become rows 268-386).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the existing row 266 (shifts 266:384 -> 268:386)
$ws.Rows("266:267").Insert()

# New row 266: Primera quality, newest date 44466 (was 44389 on the row that got pushed to 268)
$ws.Range("A266").Value = 3
$ws.Range("B266").Value = "Femacal de La Calera"
$ws.Range("C266").Value = "Coquimbo"
$ws.Range("D266").Value = 44466
$ws.Range("E266").Value = 5
$ws.Range("F266").Value = 100114014
$ws.Range("G266").Value = "Betarraga"
$ws.Range("H266").Value = "Sin especificar"
$ws.Range("I266").Value = "Primera"
$ws.Range("J266").Value = 3000
$ws.Range("K266").Value = 500
$ws.Range("L266").Value = 600
$ws.Range("M266").Value = 553
$ws.Range("N266").Value = "`$/paquete 4 unidades"
$ws.Range("O266").Value = "Provincia de Quillota"
$ws.Range("P266").Value = 138
$ws.Range("Q266").Value = 4
$ws.Range("R266").Value = "Hortaliza"

# New row 267: Segunda quality, same newest date 44466
$ws.Range("A267").Value = 3
$ws.Range("B267").Value = "Femacal de La Calera"
$ws.Range("C267").Value = "Coquimbo"
$ws.Range("D267").Value = 44466
$ws.Range("E267").Value = 5
$ws.Range("F267").Value = 100114014
$ws.Range("G267").Value = "Betarraga"
$ws.Range("H267").Value = "Sin especificar"
$ws.Range("I267").Value = "Segunda"
$ws.Range("J267").Value = 1500
$ws.Range("K267").Value = 400
$ws.Range("L267").Value = 400
$ws.Range("M267").Value = 400
$ws.Range("N267").Value = "`$/paquete 4 unidades"
$ws.Range("O267").Value = "Provincia de Quillota"
$ws.Range("P267").Value = 100
$ws.Range("Q267").Value = 4
$ws.Range("R267").Value = "Hortaliza"
